$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so values like
# "7.170" or "1.003" are not auto-converted to numbers, matching the
# original inlineStr text storage.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.971.11"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.868.06"
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "318.96"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5076"
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("D8").Value = "0.3939"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "0.08165"
$ws.Range("D10").Value = "42.12"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").Value = "22.68"
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").Value = "1.870.58"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "6.255"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "7.170"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "91.83"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").Value = "0.06409"
$ws.Range("E19").Value = "  -4.82%  "
$ws.Range("D20").Value = "17.89"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "29.971.90"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "5.807"
$ws.Range("E23").Value = "  -4.10%  "
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "2.144"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").Value = "2.089.27"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "161.13"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "20.88"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").Value = "2.229"
$ws.Range("E29").Value = "  -9.03%  "
$ws.Range("D30").Value = "127.11"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "1.058"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "0.1035"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "5.889"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").Value = "3.736"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "0.02423"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").Value = "5.212"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "0.06358"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").Value = "0.2142"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("D40").Value = "8.472"
$ws.Range("E40").Value = "  -5.91%  "
$ws.Range("D41").Value = "0.6306"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.199"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.18"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.97"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5888"
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").Value = "1.993"
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("D49").Value = "122.48"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E51").Value = "  -3.48%  "

# Restore the original (default) style on column D so no residual
# text-format styling is left behind on the cells.
$priceRange.Style = "Normal"
